$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1085
$ws1.Range("F7").Value = 1425
$ws1.Range("F8").Value = 586
$ws1.Range("F12").Value = 153
$ws1.Range("F14").Value = 426
$ws1.Range("F15").Value = 1339
$ws1.Range("F16").Value = 104
$ws1.Range("F17").Value = 96
$ws1.Range("F18").Value = 275
$ws1.Range("F20").Value = 647
$ws1.Range("F22").Value = 211
$ws1.Range("F24").Value = 5764
$ws1.Range("F26").Value = 120
$ws1.Range("F29").Value = 14363
$ws1.Range("F30").Value = 1428
$ws1.Range("F31").Value = 201
$ws1.Range("F34").Value = 4445
$ws1.Range("F35").Value = 603
$ws1.Range("F36").Value = 4190
$ws1.Range("F37").Value = 131

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1085
$ws4.Range("F7").Value = 1425
$ws4.Range("F8").Value = 586
$ws4.Range("F12").Value = 153
$ws4.Range("F14").Value = 426
$ws4.Range("F15").Value = 1339
$ws4.Range("F16").Value = 104
$ws4.Range("F17").Value = 96
$ws4.Range("F18").Value = 275
$ws4.Range("F21").Value = 647
$ws4.Range("F24").Value = 211
$ws4.Range("F27").Value = 5764
$ws4.Range("F29").Value = 120
$ws4.Range("F32").Value = 14363
$ws4.Range("F33").Value = 1428
$ws4.Range("F34").Value = 201
$ws4.Range("F37").Value = 4445
$ws4.Range("F38").Value = 603
$ws4.Range("F39").Value = 4190
$ws4.Range("F40").Value = 131

$wb.Save()
